# Update the "想去人数" (want-to-go count) values in the F column
# for both the "展览" and "全部类型" sheets, matching the regenerated
# data snapshot referenced by the commit message.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row number -> new value for column F
$updates = @{
    2  = 6632
    6  = 2025
    7  = 1540
    10 = 429
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
